$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "34.606.39"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +1.33%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.800.17"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +0.74%  "
$ws.Range("E4").Value = "  -0.30%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "226.94"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.20%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.559"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +2.05%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "33.01"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +3.73%  "
$ws.Range("E9").Value = "  +2.00%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.0695"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +1.02%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0948"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +0.27%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "2.057.08"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +0.58%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "11.13"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +0.85%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "1.770.36"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -0.89%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.639"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +2.45%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "34.576.62"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +1.39%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "4.28"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +2.52%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "68.94"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +0.98%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "248.53"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +0.79%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.0₃0801"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +3.12%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "11.40"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +4.61%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.998"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -0.26%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "4.17"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +1.61%  "
$ws.Range("E24").Value = "  +0.28%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "164.99"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +2.24%  "
$ws.Range("E26").Value = "  +1.26%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "16.56"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +1.40%  "
$ws.Range("E28").Value = "  +2.74%  "
$ws.Range("E29").Value = "  -0.31%  "
$ws.Range("E30").Value = "  +8.36%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "3.80"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +3.17%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.23"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -0.14%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.0521"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +0.38%  "
$ws.Range("E34").Value = "  +1.66%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.423.61"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -1.48%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "2.58"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +4.60%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.671"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +2.57%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.0192"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +0.61%  "
$ws.Range("E39").Value = "  +1.93%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "85.21"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +5.70%  "
$ws.Range("E41").Value = "  +1.74%  "
$ws.Range("E42").Value = "  +0.64%  "
$ws.Range("E43").Value = "  +2.43%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "13.48"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -0.85%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.0524"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +3.03%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "6.06"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -0.07%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.07"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +0.23%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.954.44"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +0.38%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "105.81"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -0.17%  "
$ws.Range("E50").Value = "  -0.34%  "
$ws.Range("E51").Value = "  -5.54%  "
